# Automatic update of files.
#
# The "Förändrad" (last-checked) date for every data row moves from
# 2026-02-17 (serial 46070) to 2026-02-19 (serial 46072).
#
# Rows 29-55 hold the same 27 "in progress" notifications as before, but the
# underlying source query re-ran and re-ordered them (a handful of brand new
# notifications appear, a few others drop their "Markägare" classification
# or pick up a different one, and some areas are corrected). We therefore
# rewrite columns A, B, F and G for rows 29-55 explicitly to match the new
# row order/content, and bump column C for every data row (2-83).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 29-55: full replacement of Beteckning/Datum/Markägare/Area (ha)
# ---------------------------------------------------------------------------
$rows2955 = @(
    @{ Row=29; A="A 57251-2025"; B=45979;               F="Övriga Aktiebolag"; G=1.6  },
    @{ Row=30; A="A 57342-2025"; B=45979;               F="Övriga Aktiebolag"; G=0.7  },
    @{ Row=31; A="A 57338-2025"; B=45979;               F="Övriga Aktiebolag"; G=0.5  },
    @{ Row=32; A="A 57246-2025"; B=45979;               F="Övriga Aktiebolag"; G=0.9  },
    @{ Row=33; A="A 57344-2025"; B=45979;               F="Övriga Aktiebolag"; G=0.7  },
    @{ Row=34; A="A 56965-2025"; B=45978;               F=$null;               G=0.6  },
    @{ Row=35; A="A 44720-2022"; B=44840;               F=$null;               G=0.9  },
    @{ Row=36; A="A 11736-2024"; B=45373;               F=$null;               G=8    },
    @{ Row=37; A="A 46451-2025"; B=45925.65525462963;   F="Holmen skog AB";    G=0.8  },
    @{ Row=38; A="A 31166-2025"; B=45832.6597337963;    F=$null;               G=3.1  },
    @{ Row=39; A="A 47659-2025"; B=45931.55524305555;   F="Holmen skog AB";    G=2.1  },
    @{ Row=40; A="A 48314-2025"; B=45933;               F=$null;               G=1    },
    @{ Row=41; A="A 32213-2025"; B=45835.58582175926;   F=$null;               G=3.5  },
    @{ Row=42; A="A 44579-2022"; B=44840.4959375;       F=$null;               G=2.8  },
    @{ Row=43; A="A 15293-2023"; B=45019;               F="Holmen skog AB";    G=0.5  },
    @{ Row=44; A="A 12447-2024"; B=45379.46974537037;   F=$null;               G=3.4  },
    @{ Row=45; A="A 17600-2022"; B=44680;               F=$null;               G=7.7  },
    @{ Row=46; A="A 59112-2025"; B=45988.48501157408;   F="Holmen skog AB";    G=1.4  },
    @{ Row=47; A="A 59222-2025"; B=45988.61856481482;   F="Holmen skog AB";    G=0.6  },
    @{ Row=48; A="A 55265-2025"; B=45968.65879629629;   F="Holmen skog AB";    G=3.5  },
    @{ Row=49; A="A 55257-2025"; B=45968.64399305556;   F=$null;               G=1.6  },
    @{ Row=50; A="A 61000-2024"; B=45645.43471064815;   F=$null;               G=10.1 },
    @{ Row=51; A="A 3337-2026";  B=46041.82646990741;   F=$null;               G=6.3  },
    @{ Row=52; A="A 56117-2022"; B=44889;               F="Holmen skog AB";    G=0.7  },
    @{ Row=53; A="A 3346-2026";  B=46041.87569444445;   F=$null;               G=3.3  },
    @{ Row=54; A="A 36709-2024"; B=45537;               F="Övriga Aktiebolag"; G=1.8  },
    @{ Row=55; A="A 36711-2024"; B=45537;               F="Övriga Aktiebolag"; G=0.7  }
)

foreach ($item in $rows2955) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A   # A: Beteckning
    $ws.Cells.Item($r, 2).Value = $item.B   # B: Datum

    if ($item.F) {
        $ws.Cells.Item($r, 6).Value = $item.F   # F: Markägare
    } else {
        $ws.Cells.Item($r, 6).ClearContents()
    }

    $ws.Cells.Item($r, 7).Value = $item.G    # G: Area (ha)
}

# ---------------------------------------------------------------------------
# 2) Every data row (2-83): bump "Förändrad" (column C) from 46070 to 46072
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 83; $r++) {
    $ws.Cells.Item($r, 3).Value = 46072
}
